$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 17424818
$ws.Range("I32").Value = 301
$ws.Range("J32").Value = 23232990
$ws.Range("K32").Value = 301
$ws.Range("L32").Value = 23232990
$ws.Range("M32").Value = 25
$ws.Range("N32").Value = -23233642
$ws.Range("H63").Value = 76238
$ws.Range("I63").Value = 21344
$ws.Range("J63").Value = 87216.8
$ws.Range("K63").Value = 21344
$ws.Range("L63").Value = 87216.8
$ws.Range("M63").Value = -20720
$ws.Range("N63").Value = -88464.8
$ws.Range("H66").Value = 76238
$ws.Range("I66").Value = 21344
$ws.Range("J66").Value = 87216.8
$ws.Range("K66").Value = 64032
$ws.Range("L66").Value = 261650.4
$ws.Range("M66").Value = -60912
$ws.Range("N66").Value = -267890.4
$ws.Range("H132").Value = 2344.8918
$ws.Range("I132").Value = 2007.7778
$ws.Range("J132").Value = 3255.1
$ws.Range("K132").Value = 6023.3334
$ws.Range("L132").Value = 9765.299999999999
$ws.Range("M132").Value = -3493.3334
$ws.Range("N132").Value = -14825.3
$ws.Range("H137").Value = 4749.5435
$ws.Range("I137").Value = 5042.6055
$ws.Range("J137").Value = 3357.5
$ws.Range("K137").Value = 15127.8165
$ws.Range("L137").Value = 10072.5
$ws.Range("M137").Value = -12577.8165
$ws.Range("N137").Value = -15172.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1145446.8
$ws.Range("I32").Value = 1285553.9
$ws.Range("J32").Value = 7076.75
$ws.Range("K32").Value = 1285553.9
$ws.Range("L32").Value = 7076.75
$ws.Range("M32").Value = -1285266.9
$ws.Range("N32").Value = -7650.75
$ws.Range("H45").Value = 999.3333
$ws.Range("I45").Value = 972.5
$ws.Range("J45").Value = 1214
$ws.Range("K45").Value = 972.5
$ws.Range("L45").Value = 1214
$ws.Range("M45").Value = -595.5
$ws.Range("N45").Value = -1968
$ws.Range("H74").Value = 277137.53
$ws.Range("I74").Value = 371871.84
$ws.Range("J74").Value = 80381.62
$ws.Range("K74").Value = 371871.84
$ws.Range("L74").Value = 80381.62
$ws.Range("M74").Value = -370997.84
$ws.Range("N74").Value = -82129.62
$ws.Range("H77").Value = 277137.53
$ws.Range("I77").Value = 371871.84
$ws.Range("J77").Value = 80381.62
$ws.Range("K77").Value = 1859359.2
$ws.Range("L77").Value = 401908.1
$ws.Range("M77").Value = -1854991.2
$ws.Range("N77").Value = -410644.1
$ws.Range("H132").Value = 31753.5
$ws.Range("I132").Value = 49158.363
$ws.Range("J132").Value = 4403
$ws.Range("K132").Value = 147475.089
$ws.Range("L132").Value = 13209
$ws.Range("M132").Value = -144945.089
$ws.Range("N132").Value = -18269

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7400
$ws.Range("I86").Value = 17533.334
$ws.Range("J86").Value = 2333.3333
$ws.Range("K86").Value = 17533.334
$ws.Range("L86").Value = 2333.3333
$ws.Range("M86").Value = -16410.334
$ws.Range("N86").Value = -4579.3333
$ws.Range("H89").Value = 7400
$ws.Range("I89").Value = 17533.334
$ws.Range("J89").Value = 2333.3333
$ws.Range("K89").Value = 87666.67
$ws.Range("L89").Value = 11666.6665
$ws.Range("M89").Value = -82050.67
$ws.Range("N89").Value = -22898.6665
$ws.Range("H94").Value = 1482.2174
$ws.Range("I94").Value = 626.6923
$ws.Range("J94").Value = 2594.4
$ws.Range("K94").Value = 626.6923
$ws.Range("L94").Value = 2594.4
$ws.Range("M94").Value = -175.6923
$ws.Range("N94").Value = -3496.4
$ws.Range("H99").Value = 6568.9
$ws.Range("I99").Value = 7076.5557
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 7076.5557
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -5578.5557
$ws.Range("N99").Value = -4996
$ws.Range("H107").Value = 1399.6666
$ws.Range("I107").Value = 1399.6666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1399.6666
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 520.3334
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 483.18182
$ws.Range("I5").Value = 125.875
$ws.Range("J5").Value = 1436
$ws.Range("K5").Value = 125.875
$ws.Range("L5").Value = 1436
$ws.Range("M5").Value = -13.875
$ws.Range("N5").Value = -1660
$ws.Range("H31").Value = 2658.544
$ws.Range("I31").Value = 1888.2
$ws.Range("J31").Value = 3884.0908
$ws.Range("K31").Value = 1888.2
$ws.Range("L31").Value = 3884.0908
$ws.Range("M31").Value = -1593.2
$ws.Range("N31").Value = -4474.0908
$ws.Range("H34").Value = 2658.544
$ws.Range("I34").Value = 1888.2
$ws.Range("J34").Value = 3884.0908
$ws.Range("K34").Value = 1888.2
$ws.Range("L34").Value = 3884.0908
$ws.Range("M34").Value = -1686.2
$ws.Range("N34").Value = -4288.0908
$ws.Range("H58").Value = 4226.528
$ws.Range("I58").Value = 4718.44
$ws.Range("J58").Value = 3108.5454
$ws.Range("K58").Value = 4718.44
$ws.Range("L58").Value = 3108.5454
$ws.Range("M58").Value = -4515.44
$ws.Range("N58").Value = -3514.5454
$ws.Range("H136").Value = 4226.528
$ws.Range("I136").Value = 4718.44
$ws.Range("J136").Value = 3108.5454
$ws.Range("K136").Value = 14155.32
$ws.Range("L136").Value = 9325.636200000001
$ws.Range("M136").Value = -11605.32
$ws.Range("N136").Value = -14425.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 9583.637000000001
$ws.Range("I132").Value = 11500
$ws.Range("J132").Value = 9157.777
$ws.Range("K132").Value = 103500
$ws.Range("L132").Value = 82419.993
$ws.Range("M132").Value = -100970
$ws.Range("N132").Value = -87479.993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3963.295
$ws.Range("I70").Value = 3739.1064
$ws.Range("J70").Value = 4303.1934
$ws.Range("K70").Value = 3739.1064
$ws.Range("L70").Value = 4303.1934
$ws.Range("M70").Value = -3469.1064
$ws.Range("N70").Value = -4843.1934
$ws.Range("H73").Value = 3963.295
$ws.Range("I73").Value = 3739.1064
$ws.Range("J73").Value = 4303.1934
$ws.Range("K73").Value = 3739.1064
$ws.Range("L73").Value = 4303.1934
$ws.Range("M73").Value = -2803.1064
$ws.Range("N73").Value = -6175.1934
$ws.Range("H132").Value = 3359.7844
$ws.Range("I132").Value = 3165.3428
$ws.Range("J132").Value = 3785.125
$ws.Range("K132").Value = 9496.028399999999
$ws.Range("L132").Value = 11355.375
$ws.Range("M132").Value = -6966.028399999999
$ws.Range("N132").Value = -16415.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1794
$ws.Range("I61").Value = 726.5
$ws.Range("J61").Value = 2648
$ws.Range("K61").Value = 726.5
$ws.Range("L61").Value = 2648
$ws.Range("M61").Value = -524.5
$ws.Range("N61").Value = -3052
$ws.Range("H113").Value = 1794
$ws.Range("I113").Value = 726.5
$ws.Range("J113").Value = 2648
$ws.Range("K113").Value = 726.5
$ws.Range("L113").Value = 2648
$ws.Range("M113").Value = 1443.5
$ws.Range("N113").Value = -6988
$ws.Range("H122").Value = 3075.2
$ws.Range("I122").Value = 3123.1428
$ws.Range("J122").Value = 2963.3333
$ws.Range("K122").Value = 9369.428400000001
$ws.Range("L122").Value = 8889.999899999999
$ws.Range("M122").Value = -6919.428400000001
$ws.Range("N122").Value = -13789.9999
$ws.Range("H132").Value = 4113.5
$ws.Range("I132").Value = 3437
$ws.Range("J132").Value = 4505.1577
$ws.Range("K132").Value = 10311
$ws.Range("L132").Value = 13515.4731
$ws.Range("M132").Value = -7781
$ws.Range("N132").Value = -18575.4731
$ws.Range("H136").Value = 3032.8542
$ws.Range("I136").Value = 2628.6562
$ws.Range("J136").Value = 3841.25
$ws.Range("K136").Value = 7885.9686
$ws.Range("L136").Value = 11523.75
$ws.Range("M136").Value = -5335.9686
$ws.Range("N136").Value = -16623.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2104.2083
$ws.Range("I81").Value = 1440.0667
$ws.Range("J81").Value = 3211.111
$ws.Range("K81").Value = 2880.1334
$ws.Range("L81").Value = 6422.222
$ws.Range("M81").Value = -1819.1334
$ws.Range("N81").Value = -8544.222
$ws.Range("H84").Value = 2104.2083
$ws.Range("I84").Value = 1440.0667
$ws.Range("J84").Value = 3211.111
$ws.Range("K84").Value = 14400.667
$ws.Range("L84").Value = 32111.11
$ws.Range("M84").Value = -9096.667000000001
$ws.Range("N84").Value = -42719.11
$ws.Range("H136").Value = 24663658
$ws.Range("I136").Value = 35755892
$ws.Range("J136").Value = 772696.4399999999
$ws.Range("K136").Value = 107267676
$ws.Range("L136").Value = 2318089.32
$ws.Range("M136").Value = -107265126
$ws.Range("N136").Value = -2323189.32
